# Update the team transition-probability matrix on Sheet1 with the refreshed
# values from games pulled March 7 (per commit message). Only the cells whose
# underlying counts changed are touched; zero/unaffected cells are left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2007299270072993
$ws.Range("C2").Value = 0.5693430656934306
$ws.Range("J2").Value = 0.01094890510948905
$ws.Range("P2").Value = 0.1313868613138686
$ws.Range("S2").Value = 0.08759124087591241
$ws.Range("B3").Value = 0.006172839506172839
$ws.Range("C3").Value = 0.04320987654320987
$ws.Range("J3").Value = 0.02469135802469136
$ws.Range("P3").Value = 0.7407407407407407
$ws.Range("S3").Value = 0.1851851851851852
$ws.Range("J4").Value = 0.05
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.07539682539682539
$ws.Range("D6").Value = 0.01587301587301587
$ws.Range("F6").Value = 0.09523809523809523
$ws.Range("J6").Value = 0.2698412698412698
$ws.Range("O6").Value = 0.02380952380952381
$ws.Range("Q6").Value = 0.1468253968253968
$ws.Range("R6").Value = 0.03571428571428571
$ws.Range("S6").Value = 0.3373015873015873
$ws.Range("B7").Value = 0.09349593495934959
$ws.Range("D7").Value = 0.02032520325203252
$ws.Range("E7").Value = 0.004065040650406504
$ws.Range("F7").Value = 0.07317073170731707
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("O7").Value = 0.008130081300813009
$ws.Range("Q7").Value = 0.2032520325203252
$ws.Range("R7").Value = 0.06504065040650407
$ws.Range("S7").Value = 0.3658536585365854
$ws.Range("B8").Value = 0.09740259740259741
$ws.Range("D8").Value = 0.01515151515151515
$ws.Range("E8").Value = 0.002164502164502165
$ws.Range("F8").Value = 0.0670995670995671
$ws.Range("J8").Value = 0.08225108225108226
$ws.Range("O8").Value = 0.02164502164502164
$ws.Range("Q8").Value = 0.2034632034632035
$ws.Range("R8").Value = 0.08658008658008658
$ws.Range("S8").Value = 0.4242424242424243
$ws.Range("B9").Value = 0.09482758620689655
$ws.Range("D9").Value = 0.008620689655172414
$ws.Range("F9").Value = 0.08620689655172414
$ws.Range("J9").Value = 0.1120689655172414
$ws.Range("O9").Value = 0.02586206896551724
$ws.Range("Q9").Value = 0.2241379310344828
$ws.Range("R9").Value = 0.1120689655172414
$ws.Range("S9").Value = 0.3362068965517241
$ws.Range("B10").Value = 0.09910714285714285
$ws.Range("D10").Value = 0.01964285714285714
$ws.Range("E10").Value = 0.002678571428571429
$ws.Range("F10").Value = 0.075
$ws.Range("J10").Value = 0.09375
$ws.Range("O10").Value = 0.01964285714285714
$ws.Range("Q10").Value = 0.2223214285714286
$ws.Range("R10").Value = 0.075
$ws.Range("S10").Value = 0.3928571428571428
$ws.Range("G11").Value = 0.1485411140583554
$ws.Range("J11").Value = 0.08753315649867374
$ws.Range("K11").Value = 0.2042440318302387
$ws.Range("L11").Value = 0.5464190981432361
$ws.Range("S11").Value = 0.01326259946949602
$ws.Range("G12").Value = 0.7546296296296297
$ws.Range("J12").Value = 0.1712962962962963
$ws.Range("K12").Value = 0.004629629629629629
$ws.Range("L12").Value = 0.04629629629629629
$ws.Range("S12").Value = 0.02314814814814815
$ws.Range("G13").Value = 0.7777777777777778
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.02222222222222222
$ws.Range("F15").Value = 0.03149606299212598
$ws.Range("H15").Value = 0.1732283464566929
$ws.Range("I15").Value = 0.09448818897637795
$ws.Range("J15").Value = 0.3031496062992126
$ws.Range("K15").Value = 0.09448818897637795
$ws.Range("M15").Value = 0.007874015748031496
$ws.Range("O15").Value = 0.04330708661417323
$ws.Range("S15").Value = 0.2519685039370079
$ws.Range("F16").Value = 0.01685393258426966
$ws.Range("H16").Value = 0.2022471910112359
$ws.Range("I16").Value = 0.07303370786516854
$ws.Range("J16").Value = 0.4213483146067415
$ws.Range("K16").Value = 0.09550561797752809
$ws.Range("M16").Value = 0.01685393258426966
$ws.Range("O16").Value = 0.07303370786516854
$ws.Range("S16").Value = 0.101123595505618
$ws.Range("F17").Value = 0.02892561983471074
$ws.Range("H17").Value = 0.2148760330578512
$ws.Range("I17").Value = 0.09090909090909091
$ws.Range("J17").Value = 0.3223140495867768
$ws.Range("K17").Value = 0.1198347107438017
$ws.Range("M17").Value = 0.02892561983471074
$ws.Range("N17").Value = 0.004132231404958678
$ws.Range("O17").Value = 0.08057851239669421
$ws.Range("S17").Value = 0.109504132231405
$ws.Range("F18").Value = 0.005714285714285714
$ws.Range("H18").Value = 0.1771428571428571
$ws.Range("I18").Value = 0.12
$ws.Range("J18").Value = 0.32
$ws.Range("K18").Value = 0.1371428571428571
$ws.Range("M18").Value = 0.005714285714285714
$ws.Range("O18").Value = 0.08
$ws.Range("S18").Value = 0.1542857142857143
$ws.Range("F19").Value = 0.02129337539432177
$ws.Range("H19").Value = 0.194006309148265
$ws.Range("I19").Value = 0.1025236593059937
$ws.Range("J19").Value = 0.3170347003154574
$ws.Range("K19").Value = 0.1324921135646688
$ws.Range("M19").Value = 0.022870662460567823
$ws.Range("N19").Value = 0.003154574132492113
$ws.Range("O19").Value = 0.07413249211356467
$ws.Range("S19").Value = 0.1324921135646688
